# Fruta / hortaliza, semanal
# Insert a new data row at row 299 (pushing the existing rows 299-372 down
# to 300-373) and populate it with the new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(299).Insert()

$ws.Range("A299").Value = 10
$ws.Range("B299").Value = "Vega Modelo de Temuco"
$ws.Range("C299").Value = "La Araucanía"
$ws.Range("D299").Value = 44736
$ws.Range("E299").Value = 9
$ws.Range("F299").Value = "Fruta"
$ws.Range("G299").Value = 100108
$ws.Range("H299").Value = "Tropicales y subtropicales"
$ws.Range("I299").Value = 100108002
$ws.Range("J299").Value = "Mango"
$ws.Range("K299").Value = "Sin especificar"
$ws.Range("L299").Value = "Primera"
$ws.Range("M299").Value = 800
$ws.Range("N299").Value = 9000
$ws.Range("O299").Value = 9000
$ws.Range("P299").Value = 9000
$ws.Range("Q299").Value = "$/bandeja 4 kilos"
$ws.Range("R299").Value = "Brasil"
$ws.Range("S299").Value = 2250
$ws.Range("T299").Value = 4
